$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds metadata for "Tipo de estudios realizados".
# C3: reclassified from a "measure" to a "dimension" in the iaest vocabulary.
$ws.Range("C3").Value = "iaest-dimension:tipo-de-estudios-realizados"

# C4: "medida" (measure) -> "dim" (dimension), matching the new classification.
$ws.Range("C4").Value = "dim"

# C5: datatype changes from a plain xsd:string to a coded skos:Concept.
$ws.Range("C5").Value = "skos:Concept"

# New row 6: reference the external code-list mapping file for this dimension.
$ws.Range("C6").Value = "mapping-tipo-de-estudios-realizados.xlsx"

# Match the formatting used by the rest of the metadata column (style used by C5)
# instead of the workbook's bare default style.
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)
